# AO2, atividade encontro remoto 2 e encontro remoto 3
#
# Remove the unused trailing "responsável" template slide (sldId 410 /
# Lorem-ipsum placeholder content) — slide 7 of 7 — along with its
# notes page. This mirrors deleting the slide from the Slides pane in
# PowerPoint, which also drops its entry from the slide-id list and its
# parts (slide7.xml, notesSlide7.xml) from the package.
$p = $ppt.ActivePresentation
$p.Slides.Item(7).Delete()
